# Remove the "checkDates" method row from Sheet 2 (Class: Customer).
# This deletes the entire row, shifting subsequent rows up by one,
# which matches the target edit (row "checkDates" removed, all rows
# below it shift from row 9-19 to row 8-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 2 - Class_ Customer - Cla")
$ws.Rows.Item(8).Delete()
